$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for Alice Johnson / CS2025_00073 (duplicate "Login Issues." record)
# so all subsequent rows shift up by one.
$ws.Range("A24").EntireRow.Delete()

# Renumber CaseID column sequentially (CS2025_00142 .. CS2025_00171) and
# refresh the "Ai Generated Sub Categories" column with the cleaned-up /
# retrained category labels.
$caseIds = @(
    "CS2025_00142","CS2025_00143","CS2025_00144","CS2025_00145","CS2025_00146",
    "CS2025_00147","CS2025_00148","CS2025_00149","CS2025_00150","CS2025_00151",
    "CS2025_00152","CS2025_00153","CS2025_00154","CS2025_00155","CS2025_00156",
    "CS2025_00157","CS2025_00158","CS2025_00159","CS2025_00160","CS2025_00161",
    "CS2025_00162","CS2025_00163","CS2025_00164","CS2025_00165","CS2025_00166",
    "CS2025_00167","CS2025_00168","CS2025_00169","CS2025_00170","CS2025_00171"
)

$subCategories = @(
    "Authentication Issue","Software Troubleshooting.","Refund Request","Payment Disputes","Product Features.",
    "Product Pricing.","Hardware Issue.","Website Maintenance","Invoice Request.","Operating Hours",
    "App Troubleshooting.","Refund Request","Billing Inquiry.","Network Troubleshooting.","Account Management",
    "Software Installation.","Store Location.","Email Setup.","Payment Plans.","Password Reset",
    "Product Comparison.","Payment Issue","Hardware Issue.","Product Support.","Performance Issue.",
    "Account Management","Training Services","Account Access.","Account Inquiries","Technical Support."
)

for ($i = 0; $i -lt $caseIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $caseIds[$i]
    $ws.Cells.Item($row, 7).Value = $subCategories[$i]
}
